# Automatische test-sync: 2025-08-19 20:44:50
# Adds a new log row (row 18) to the "Logs" sheet, extends the conditional
# formatting ranges that cover the data rows, and refreshes the aggregated
# count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the "Logs" sheet -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A18").Value = "Interne taak"
$logs.Range("B18").Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("D18").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F18").Value = "2025-08-19 20:44:25"
$logs.Range("G18").Value = "Nee"
$logs.Range("H18").Value = "Ja"
$logs.Range("I18").Value = "Nee"
$logs.Range("J18").Value = "Nee"

# --- 2. Extend the conditional formatting ranges from row 17 to row 18 -----
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "17")
    $newRange = $logs.Range($col + "2:" + $col + "18")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the aggregated total on the "Dashboard" sheet ---------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 17
